$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "59.350.73"
Set-TextValue $ws.Range("E2") "  -1.72%  "
Set-TextValue $ws.Range("D3") "2.573.19"
Set-TextValue $ws.Range("E3") "  -2.12%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.27%  "
Set-TextValue $ws.Range("D5") "554.28"
Set-TextValue $ws.Range("E5") "  -2.62%  "
Set-TextValue $ws.Range("D6") "141.62"
Set-TextValue $ws.Range("E6") "  -2.96%  "
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("E8") "  -1.44%  "
Set-TextValue $ws.Range("D9") "2.575.47"
Set-TextValue $ws.Range("E9") "  -2.81%  "
Set-TextValue $ws.Range("D10") "6.73"
Set-TextValue $ws.Range("E10") "  -1.57%  "
Set-TextValue $ws.Range("E11") "  -1.30%  "
Set-TextValue $ws.Range("D12") "0.164"
Set-TextValue $ws.Range("E12") "  +9.70%  "
Set-TextValue $ws.Range("D13") "0.352"
Set-TextValue $ws.Range("E13") "  +1.87%  "
Set-TextValue $ws.Range("D14") "3.025.26"
Set-TextValue $ws.Range("E14") "  -2.19%  "
Set-TextValue $ws.Range("D15") "59.325.57"
Set-TextValue $ws.Range("E15") "  -1.76%  "
Set-TextValue $ws.Range("D16") "23.00"
Set-TextValue $ws.Range("E16") "  +4.03%  "
Set-TextValue $ws.Range("E17") "  -0.18%  "
Set-TextValue $ws.Range("D18") "2.573.29"
Set-TextValue $ws.Range("E18") "  -2.61%  "
Set-TextValue $ws.Range("E19") "  +0.04%  "
Set-TextValue $ws.Range("D20") "338.13"
Set-TextValue $ws.Range("E20") "  -1.00%  "
Set-TextValue $ws.Range("E21") "  -1.02%  "
Set-TextValue $ws.Range("D22") "6.45"
Set-TextValue $ws.Range("E22") "  +1.29%  "
Set-TextValue $ws.Range("D23") "0.997"
Set-TextValue $ws.Range("E23") "  -0.24%  "
Set-TextValue $ws.Range("D24") "0.478"
Set-TextValue $ws.Range("E24") "  +7.64%  "
Set-TextValue $ws.Range("D25") "62.67"
Set-TextValue $ws.Range("E25") "  -5.07%  "
Set-TextValue $ws.Range("D26") "0.998"
Set-TextValue $ws.Range("E26") "  -0.29%  "
Set-TextValue $ws.Range("E27") "  -2.88%  "
Set-TextValue $ws.Range("D28") "7.38"
Set-TextValue $ws.Range("E28") "  +0.11%  "
Set-TextValue $ws.Range("D29") "0.0₃0771"
Set-TextValue $ws.Range("E29") "  -3.94%  "
Set-TextValue $ws.Range("E30") "  +0.02%  "
Set-TextValue $ws.Range("D31") "6.18"
Set-TextValue $ws.Range("E31") "  +0.85%  "
Set-TextValue $ws.Range("E32") "  -2.50%  "
Set-TextValue $ws.Range("D33") "158.70"
Set-TextValue $ws.Range("E33") "  +0.20%  "
Set-TextValue $ws.Range("D34") "19.03"
Set-TextValue $ws.Range("E34") "  -0.58%  "
Set-TextValue $ws.Range("E35") "  -0.50%  "
Set-TextValue $ws.Range("D36") "1.17"
Set-TextValue $ws.Range("E36") "  +1.25%  "
Set-TextValue $ws.Range("E37") "  +0.83%  "
Set-TextValue $ws.Range("D38") "37.40"
Set-TextValue $ws.Range("E38") "  -0.27%  "
Set-TextValue $ws.Range("D39") "0.854"
Set-TextValue $ws.Range("E39") "  -2.91%  "
Set-TextValue $ws.Range("D40") "1.47"
Set-TextValue $ws.Range("E40") "  -2.96%  "
Set-TextValue $ws.Range("E41") "  +0.58%  "
Set-TextValue $ws.Range("D42") "289.45"
Set-TextValue $ws.Range("E42") "  -3.65%  "
Set-TextValue $ws.Range("D43") "137.80"
Set-TextValue $ws.Range("E43") "  +8.24%  "
Set-TextValue $ws.Range("D44") "0.998"
Set-TextValue $ws.Range("E44") "  +0.35%  "
Set-TextValue $ws.Range("E45") "  -1.38%  "
Set-TextValue $ws.Range("D46") "0.591"
Set-TextValue $ws.Range("E46") "  -1.94%  "
Set-TextValue $ws.Range("E47") "  -0.24%  "
Set-TextValue $ws.Range("D48") "0.0529"
Set-TextValue $ws.Range("E48") "  -2.73%  "
Set-TextValue $ws.Range("E49") "  -1.04%  "
Set-TextValue $ws.Range("D50") "18.62"
Set-TextValue $ws.Range("E50") "  -0.49%  "
Set-TextValue $ws.Range("D51") "1.950.98"
Set-TextValue $ws.Range("E51") "  -0.72%  "
